$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.275.10'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.07%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.636.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.07%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.04%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.13'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.28%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('E8').Value = '  -1.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.633.79'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.170'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.96%  '

$ws.Range('E11').Value = '  +0.98%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.39%  '

$ws.Range('E13').Value = '  -2.19%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.116.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.08%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '72.136.77'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.13%  '

$ws.Range('E16').Value = '  -0.91%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.79'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.75%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.641.86'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.82%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.42'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.97%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.72%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '374.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.92%  '

$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.70'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.19%  '

$ws.Range('E26').Value = '  -2.72%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.54'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.45%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.770.40'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.22%  '

$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('E30').Value = '  +0.43%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '496.25'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.21%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.27'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.79'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.91%  '

$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.16%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.18'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.28%  '

$ws.Range('E38').Value = '  +3.41%  '

$ws.Range('E39').Value = '  -1.51%  '

$ws.Range('E40').Value = '  -2.38%  '

$ws.Range('E41').Value = '  +0.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.72'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.54%  '

$ws.Range('E43').Value = '  -1.21%  '

$ws.Range('E44').Value = '  -2.86%  '

$ws.Range('E45').Value = '  -2.40%  '

$ws.Range('E46').Value = '  -0.51%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '151.99'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.19%  '

$ws.Range('E48').Value = '  -2.32%  '

$ws.Range('E49').Value = '  -0.55%  '

$ws.Range('E50').Value = '  -2.23%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.599'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.37%  '
